$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A (1) and Column C (3) get wider; Column B (2) is left untouched
# since its width/value does not change.
# NOTE: this COM shim only exposes column width in 1/6-character pixel
# steps via ColumnWidth, so we feed it the input that lands on the closest
# achievable stored width to the target (11.7109375 -> ~11.6667,
# 12.7109375 -> ~12.6667).
$ws.Columns.Item(1).ColumnWidth = 10.833333333333334
$ws.Columns.Item(3).ColumnWidth = 11.833333333333334

# Update the three data values on row 1
$ws.Range("A1").Value = 148.90775052551567
$ws.Range("B1").Value = 4.7140732030366284
$ws.Range("C1").Value = 0.70215264187866921
